$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at J (this naturally shifts the old, data-less
# width-only column definitions that used to be at K/N to L/O)
$ws.Columns("J:J").Insert()

# Header for the new "Populate" column
$ws.Range("J1").Value = "Populate"

# Per-row Populate flag: everything is populated ("Y") except the two
# "do not stuff" rows (row 8 = U3 custom QFN reference only, row 18 = U5
# SCuM QFN custom part), which are marked "N"
$populateByRow = @{
    2  = "Y"; 3  = "Y"; 4  = "Y"; 5  = "Y"; 6  = "Y"; 7  = "Y";
    8  = "N";
    9  = "Y"; 10 = "Y"; 11 = "Y"; 12 = "Y"; 13 = "Y"; 14 = "Y"; 15 = "Y";
    16 = "Y"; 17 = "Y";
    18 = "N";
    19 = "Y"; 20 = "Y"; 21 = "Y"; 22 = "Y"; 23 = "Y"; 24 = "Y"; 25 = "Y"
}

foreach ($row in 2..25) {
    $ws.Range("J$row").Value = $populateByRow[$row]
}

# Rows 8 and 18 are the "do not populate" rows - highlight the whole row
# (only the cells that actually hold data) with Excel's built-in
# "Neutral" cell style (yellow fill / brown text)
$noPopulateCols = @{
    "8"  = @("A", "B", "C", "E", "F", "H", "I", "J");
    "18" = @("A", "C", "E", "I", "J")
}
foreach ($row in 8, 18) {
    foreach ($col in $noPopulateCols["$row"]) {
        $ws.Range("$col$row").Style = "Neutral"
    }
}

# Roughly re-fit the QTY/BOARD column now that a new column sits next to it
$ws.Range("I1:I25").EntireColumn.AutoFit()

# Leave the cursor where the author left it after editing
$ws.Range("D35").Select()
